$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 9 (Leve Item ID 5487)
$ws.Cells.Item(9, 8).Value = 525.1177
$ws.Cells.Item(9, 9).Value = 395.57144
$ws.Cells.Item(9, 11).Value = 395.57144
$ws.Cells.Item(9, 13).Value = -226.57144
# Row 32 (Leve Item ID 5484)
$ws.Cells.Item(32, 8).Value = 6078.7144
$ws.Cells.Item(32, 9).Value = 2710.4
$ws.Cells.Item(32, 10).Value = 14499.5
$ws.Cells.Item(32, 11).Value = 2710.4
$ws.Cells.Item(32, 12).Value = 14499.5
$ws.Cells.Item(32, 13).Value = -2384.4
$ws.Cells.Item(32, 14).Value = -15151.5
# Row 40 (Leve Item ID 5505)
$ws.Cells.Item(40, 8).Value = 6053.1035
$ws.Cells.Item(40, 9).Value = 4599.3335
$ws.Cells.Item(40, 10).Value = 8432
$ws.Cells.Item(40, 11).Value = 4599.3335
$ws.Cells.Item(40, 12).Value = 8432
$ws.Cells.Item(40, 13).Value = -4424.3335
$ws.Cells.Item(40, 14).Value = -8782
# Row 43 (Leve Item ID 5472)
$ws.Cells.Item(43, 8).Value = 6753.731
$ws.Cells.Item(43, 10).Value = 4759.4
$ws.Cells.Item(43, 12).Value = 4759.4
$ws.Cells.Item(43, 14).Value = -4897.4
# Row 51 (Leve Item ID 5486)
$ws.Cells.Item(51, 8).Value = 12604.235
$ws.Cells.Item(51, 9).Value = 15999.333
$ws.Cells.Item(51, 11).Value = 15999.333
$ws.Cells.Item(51, 13).Value = -15515.333
# Row 53 (Leve Item ID 5479)
$ws.Cells.Item(53, 8).Value = 393.2857
$ws.Cells.Item(53, 9).Value = 290.22223
$ws.Cells.Item(53, 10).Value = 578.8
$ws.Cells.Item(53, 11).Value = 290.22223
$ws.Cells.Item(53, 12).Value = 578.8
$ws.Cells.Item(53, 13).Value = 346.77777
$ws.Cells.Item(53, 14).Value = -1852.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5 (Leve Item ID 5091)
$ws.Cells.Item(5, 8).Value = 167.4
$ws.Cells.Item(5, 9).Value = 169.5
$ws.Cells.Item(5, 11).Value = 169.5
$ws.Cells.Item(5, 13).Value = -57.5
# Row 74 (Leve Item ID 44000)
$ws.Cells.Item(74, 8).Value = 4888.36
$ws.Cells.Item(74, 9).Value = 3104.625
$ws.Cells.Item(74, 11).Value = 3104.625
$ws.Cells.Item(74, 13).Value = -2230.625
# Row 77 (Leve Item ID 44000)
$ws.Cells.Item(77, 8).Value = 4888.36
$ws.Cells.Item(77, 9).Value = 3104.625
$ws.Cells.Item(77, 11).Value = 15523.125
$ws.Cells.Item(77, 13).Value = -11155.125
# Row 132 (Leve Item ID 43997)
$ws.Cells.Item(132, 8).Value = 3434.9473
$ws.Cells.Item(132, 9).Value = 2356.0625
$ws.Cells.Item(132, 10).Value = 9189
$ws.Cells.Item(132, 11).Value = 7068.1875
$ws.Cells.Item(132, 12).Value = 27567
$ws.Cells.Item(132, 13).Value = -4538.1875
$ws.Cells.Item(132, 14).Value = -32627

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4 (Leve Item ID 5091)
$ws.Cells.Item(4, 8).Value = 167.4
$ws.Cells.Item(4, 9).Value = 169.5
$ws.Cells.Item(4, 11).Value = 169.5
$ws.Cells.Item(4, 13).Value = -54.5
# Row 22 (Leve Item ID 5092)
$ws.Cells.Item(22, 8).Value = 1260.6428
$ws.Cells.Item(22, 9).Value = 1084.3334
$ws.Cells.Item(22, 10).Value = 1578
$ws.Cells.Item(22, 11).Value = 1084.3334
$ws.Cells.Item(22, 12).Value = 1578
$ws.Cells.Item(22, 13).Value = -911.3334
$ws.Cells.Item(22, 14).Value = -1924
# Row 26 (Leve Item ID 19535)
$ws.Cells.Item(26, 8).Value = 60000
$ws.Cells.Item(26, 9).Value = 60000
$ws.Cells.Item(26, 11).Value = 60000
$ws.Cells.Item(26, 13).Value = -59708
# Row 134 (Leve Item ID 43998)
$ws.Cells.Item(134, 8).Value = 3672.4048
$ws.Cells.Item(134, 9).Value = 3305.2974
$ws.Cells.Item(134, 11).Value = 9915.8922
$ws.Cells.Item(134, 13).Value = -7380.8922

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Cells.Item(7, 8).Value = 865.7083
$ws.Cells.Item(7, 9).Value = 845.82355
$ws.Cells.Item(7, 11).Value = 845.82355
$ws.Cells.Item(7, 13).Value = -732.82355
# Row 14 (Leve Item ID 1998)
$ws.Cells.Item(14, 8).Value = 300
$ws.Cells.Item(14, 10).Value = 300
$ws.Cells.Item(14, 12).Value = 300
$ws.Cells.Item(14, 14).Value = -640
# Row 58 (Leve Item ID 44021)
$ws.Cells.Item(58, 8).Value = 2430.889
$ws.Cells.Item(58, 9).Value = 1770.8948
$ws.Cells.Item(58, 11).Value = 1770.8948
$ws.Cells.Item(58, 13).Value = -1567.8948
# Row 99 (Leve Item ID 36198)
$ws.Cells.Item(99, 8).Value = 4492.722
$ws.Cells.Item(99, 9).Value = 4479.8
$ws.Cells.Item(99, 10).Value = 4508.875
$ws.Cells.Item(99, 11).Value = 4479.8
$ws.Cells.Item(99, 12).Value = 4508.875
$ws.Cells.Item(99, 13).Value = -2981.8
$ws.Cells.Item(99, 14).Value = -7504.875
# Row 105 (Leve Item ID 19928)
$ws.Cells.Item(105, 8).Value = 1936
$ws.Cells.Item(105, 9).Value = 2096.182
$ws.Cells.Item(105, 11).Value = 2096.182
$ws.Cells.Item(105, 13).Value = -349.1819999999998
# Row 126 (Leve Item ID 36198)
$ws.Cells.Item(126, 8).Value = 4492.722
$ws.Cells.Item(126, 9).Value = 4479.8
$ws.Cells.Item(126, 10).Value = 4508.875
$ws.Cells.Item(126, 11).Value = 13439.4
$ws.Cells.Item(126, 12).Value = 13526.625
$ws.Cells.Item(126, 13).Value = -10969.4
$ws.Cells.Item(126, 14).Value = -18466.625
# Row 136 (Leve Item ID 44021)
$ws.Cells.Item(136, 8).Value = 2430.889
$ws.Cells.Item(136, 9).Value = 1770.8948
$ws.Cells.Item(136, 11).Value = 5312.6844
$ws.Cells.Item(136, 13).Value = -2762.6844
# Row 141 (Leve Item ID 43345)
$ws.Cells.Item(141, 8).Value = 259988.78
$ws.Cells.Item(141, 10).Value = 341264.7
$ws.Cells.Item(141, 12).Value = 341264.7
$ws.Cells.Item(141, 14).Value = -351624.7

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 17 (Leve Item ID 4640)
$ws.Cells.Item(17, 8).Value = 1059.0834
$ws.Cells.Item(17, 9).Value = 1101
$ws.Cells.Item(17, 10).Value = 933.3333
$ws.Cells.Item(17, 11).Value = 3303
$ws.Cells.Item(17, 12).Value = 2799.9999
$ws.Cells.Item(17, 13).Value = -3134
$ws.Cells.Item(17, 14).Value = -3137.9999
# Row 33 (Leve Item ID 4867)
$ws.Cells.Item(33, 8).Value = 37
$ws.Cells.Item(33, 9).Value = 15
$ws.Cells.Item(33, 10).Value = 40.666668
$ws.Cells.Item(33, 11).Value = 90
$ws.Cells.Item(33, 12).Value = 244.000008
$ws.Cells.Item(33, 13).Value = 193
$ws.Cells.Item(33, 14).Value = -810.000008
# Row 63 (Leve Item ID 12866)
$ws.Cells.Item(63, 8).Value = 50000
$ws.Cells.Item(63, 9).Value = 50000
$ws.Cells.Item(63, 11).Value = 150000
$ws.Cells.Item(63, 13).Value = -149251
# Row 66 (Leve Item ID 12866)
$ws.Cells.Item(66, 8).Value = 50000
$ws.Cells.Item(66, 9).Value = 50000
$ws.Cells.Item(66, 11).Value = 450000
$ws.Cells.Item(66, 13).Value = -446256
# Row 121 (Leve Item ID 27878)
$ws.Cells.Item(121, 8).Value = 1268138.9
$ws.Cells.Item(121, 9).Value = 3906
$ws.Cells.Item(121, 10).Value = 2110960.8
$ws.Cells.Item(121, 11).Value = 11718
$ws.Cells.Item(121, 12).Value = 6332882.399999999
$ws.Cells.Item(121, 13).Value = -10408
$ws.Cells.Item(121, 14).Value = -6335502.399999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Cells.Item(70, 8).Value = 5241.933
$ws.Cells.Item(70, 9).Value = 5013.1
$ws.Cells.Item(70, 10).Value = 5699.6
$ws.Cells.Item(70, 11).Value = 5013.1
$ws.Cells.Item(70, 12).Value = 5699.6
$ws.Cells.Item(70, 13).Value = -4743.1
$ws.Cells.Item(70, 14).Value = -6239.6
# Row 73 (Leve Item ID 14146)
$ws.Cells.Item(73, 8).Value = 5241.933
$ws.Cells.Item(73, 9).Value = 5013.1
$ws.Cells.Item(73, 10).Value = 5699.6
$ws.Cells.Item(73, 11).Value = 5013.1
$ws.Cells.Item(73, 12).Value = 5699.6
$ws.Cells.Item(73, 13).Value = -4077.1
$ws.Cells.Item(73, 14).Value = -7571.6
# Row 80 (Leve Item ID 12521)
$ws.Cells.Item(80, 8).Value = 2449.125
$ws.Cells.Item(80, 10).Value = 2914.3333
$ws.Cells.Item(80, 12).Value = 2914.3333
$ws.Cells.Item(80, 14).Value = -4910.3333
# Row 83 (Leve Item ID 12521)
$ws.Cells.Item(83, 8).Value = 2449.125
$ws.Cells.Item(83, 10).Value = 2914.3333
$ws.Cells.Item(83, 12).Value = 14571.6665
$ws.Cells.Item(83, 14).Value = -24555.6665
# Row 113 (Leve Item ID 27710)
$ws.Cells.Item(113, 8).Value = 4359.875
$ws.Cells.Item(113, 9).Value = 2794.75
$ws.Cells.Item(113, 11).Value = 2794.75
$ws.Cells.Item(113, 13).Value = -624.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 13 (Leve Item ID 3546)
$ws.Cells.Item(13, 8).Value = 11328.667
# Row 16 (Leve Item ID 5289)
$ws.Cells.Item(16, 8).Value = 1416.4615
$ws.Cells.Item(16, 9).Value = 1416.4615
$ws.Cells.Item(16, 11).Value = 1416.4615
$ws.Cells.Item(16, 13).Value = -1246.4615
# Row 22 (Leve Item ID 5277)
$ws.Cells.Item(22, 8).Value = 2445.7273
$ws.Cells.Item(22, 9).Value = 2701.1428
$ws.Cells.Item(22, 10).Value = 1998.75
$ws.Cells.Item(22, 11).Value = 2701.1428
$ws.Cells.Item(22, 12).Value = 1998.75
$ws.Cells.Item(22, 13).Value = -2406.1428
$ws.Cells.Item(22, 14).Value = -2588.75
# Row 27 (Leve Item ID 5277)
$ws.Cells.Item(27, 8).Value = 2445.7273
$ws.Cells.Item(27, 9).Value = 2701.1428
$ws.Cells.Item(27, 10).Value = 1998.75
$ws.Cells.Item(27, 11).Value = 2701.1428
$ws.Cells.Item(27, 12).Value = 1998.75
$ws.Cells.Item(27, 13).Value = -2594.1428
$ws.Cells.Item(27, 14).Value = -2212.75
# Row 55 (Leve Item ID 5284)
$ws.Cells.Item(55, 8).Value = 2105.4443
$ws.Cells.Item(55, 10).Value = 2933.3333
$ws.Cells.Item(55, 12).Value = 2933.3333
$ws.Cells.Item(55, 14).Value = -3279.3333
# Row 82 (Leve Item ID 12565)
$ws.Cells.Item(82, 8).Value = 2560.6365
$ws.Cells.Item(82, 9).Value = 1575.5
$ws.Cells.Item(82, 11).Value = 1575.5
$ws.Cells.Item(82, 13).Value = -1214.5
# Row 85 (Leve Item ID 12565)
$ws.Cells.Item(85, 8).Value = 2560.6365
$ws.Cells.Item(85, 9).Value = 1575.5
$ws.Cells.Item(85, 11).Value = 1575.5
$ws.Cells.Item(85, 13).Value = -327.5
# Row 122 (Leve Item ID 36247)
$ws.Cells.Item(122, 8).Value = 5342.077
$ws.Cells.Item(122, 9).Value = 4518.263
$ws.Cells.Item(122, 11).Value = 13554.789
$ws.Cells.Item(122, 13).Value = -11104.789
# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 8306.368
$ws.Cells.Item(132, 9).Value = 7227.7617
$ws.Cells.Item(132, 10).Value = 9638.764999999999
$ws.Cells.Item(132, 11).Value = 21683.2851
$ws.Cells.Item(132, 12).Value = 28916.295
$ws.Cells.Item(132, 13).Value = -19153.2851
$ws.Cells.Item(132, 14).Value = -33976.295
# Row 136 (Leve Item ID 44060)
$ws.Cells.Item(136, 8).Value = 7516.163
$ws.Cells.Item(136, 9).Value = 5268.6816
$ws.Cells.Item(136, 10).Value = 9347.444
$ws.Cells.Item(136, 11).Value = 15806.0448
$ws.Cells.Item(136, 12).Value = 28042.332
$ws.Cells.Item(136, 13).Value = -13256.0448
$ws.Cells.Item(136, 14).Value = -33142.33199999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Cells.Item(132, 8).Value = 3305.04
$ws.Cells.Item(132, 9).Value = 2769.3914
$ws.Cells.Item(132, 11).Value = 8308.174199999999
$ws.Cells.Item(132, 13).Value = -5778.174199999999
# Row 136 (Leve Item ID 44031)
$ws.Cells.Item(136, 8).Value = 3961.8718
$ws.Cells.Item(136, 9).Value = 4203.136
$ws.Cells.Item(136, 11).Value = 12609.408
$ws.Cells.Item(136, 13).Value = -10059.408
